# Apply the edits described by the diff:
#  - Summary sheet: Strategy Total PnL (B4) 100 -> 50
#  - Symbols sheet: AAPL -> SPY row updated (pnl 100 -> 50, description updated)
#  - Strategies sheet: AAPL -> SPY row updated (expiry, strategy, pnl, hold_days, theta_per_day, description)

$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B4").Value = 50

# --- Symbols sheet ---
$symbols = $wb.Worksheets.Item("Symbols")
$symbols.Range("A2").Value = "SPY"
$symbols.Range("B2").Value = 50
$symbols.Range("E2").Value = "Options on S&P 500 ETF"

# --- Strategies sheet ---
$strategies = $wb.Worksheets.Item("Strategies")
$strategies.Range("A2").Value = "SPY"

# B2 holds a text date-like string ("2025-01-10"). A plain .Value assignment
# lets Excel auto-detect it as a real date and reformat the cell, so force
# text formatting for the write and then clear the format override back to
# the sheet default so no stray style is left behind.
$strategies.Range("B2").NumberFormat = "@"
$strategies.Range("B2").Value = "2025-01-10"
$strategies.Range("B2").ClearFormats()

$strategies.Range("C2").Value = "Short Put"
$strategies.Range("D2").Value = 50
$strategies.Range("E2").Value = 0.08333333333333333
$strategies.Range("F2").Value = 600
$strategies.Range("G2").Value = "Options on S&P 500 ETF"
